$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# New row 53: 2012-12-04 (serial 41247), 2.5h effort, new description
$ws.Cells.Item(53, 1).Value = 41247
$ws.Cells.Item(53, 1).NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Cells.Item(53, 2).Value = 2.5
$ws.Cells.Item(53, 4).Value = "Manual: new section continued, new figure for illustration"

# New row 54: 2012-12-05 (serial 41248), 0.75h effort, 2.5h additional effort, new description
$ws.Cells.Item(54, 1).Value = 41248
$ws.Cells.Item(54, 1).NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Cells.Item(54, 2).Value = 0.75
$ws.Cells.Item(54, 3).Value = 2.5
$ws.Cells.Item(54, 4).Value = "Manual: new section 2.7 completed, including new figures and first review"

# New row 55: 2012-12-06 (serial 41249), 2.5h effort, new description
$ws.Cells.Item(55, 1).Value = 41249
$ws.Cells.Item(55, 1).NumberFormat = "ddd\ dd/mm/yyyy"
$ws.Cells.Item(55, 2).Value = 2.5
$ws.Cells.Item(55, 4).Value = "Manual: Section 4.4, data type system time rewritten"

$ws.Range("D55").Select()
